# Estadisticos Matutinos 15 Oct
# Fills in BIOLOGIA grades that were pending (-1), recomputes the BIOLOGIA /
# ETICA rows on "Totales", swaps the BIOLOGIA <-> ETICA teacher assignment
# wherever it shows up on "Blancos", and refreshes "Rescatables" now that
# several students are no longer missing a grade.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Calificaciones: BIOLOGIA grades (column D / 1P, column V / Final) that
#    were placeholders (-1) now have a real grade.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Calificaciones")

$bioGrades = @{
    11 = 9
    13 = 7
    16 = 7
    17 = 8
    19 = 6
    21 = 7
    22 = 6
    23 = 10
    24 = 10
}

foreach ($r in $bioGrades.Keys) {
    $g = $bioGrades[$r]
    $ws1.Cells.Item($r, 4).Value = $g   # column D
    $ws1.Cells.Item($r, 22).Value = $g  # column V
}

# ---------------------------------------------------------------------
# 2) Totales: ETICA (row 2) and BIOLOGIA (row 3) summary rows recomputed
#    now that BIOLOGIA has real grades.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Totales")

# Row 2 -> ETICA
$ws2.Cells.Item(2, 1).Value = "ÉTICA"
$ws2.Cells.Item(2, 3).Value = 23
$ws2.Cells.Item(2, 4).Value = 8
$ws2.Cells.Item(2, 5).Value = 0
$ws2.Cells.Item(2, 6).Value = 34.78
$ws2.Cells.Item(2, 7).Value = 0
$ws2.Cells.Item(2, 8).Value = 8.9
$ws2.Cells.Item(2, 9).Value = 15
$ws2.Cells.Item(2, 10).Value = 65.22

# Row 3 -> BIOLOGIA
$ws2.Cells.Item(3, 1).Value = "BIOLOGÍA"
$ws2.Cells.Item(3, 3).Value = 24
$ws2.Cells.Item(3, 4).Value = 9
$ws2.Cells.Item(3, 5).Value = 0
$ws2.Cells.Item(3, 6).Value = 37.5
$ws2.Cells.Item(3, 7).Value = 0
$ws2.Cells.Item(3, 8).Value = 7.8
$ws2.Cells.Item(3, 9).Value = 15
$ws2.Cells.Item(3, 10).Value = 62.5

# ---------------------------------------------------------------------
# 3) Blancos: teacher assigned to BIOLOGIA / ETICA swaps on every row of
#    the sheet (column E = Materia, column F = Docente).
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Blancos")

$lastRow3 = $ws3.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow3; $r++) {
    $materia = $ws3.Cells.Item($r, 5).Text
    if ($materia -eq "BIOLOGÍA") {
        $ws3.Cells.Item($r, 6).Value = "Castro Vasquez Julieta"
    } elseif ($materia -eq "ÉTICA") {
        $ws3.Cells.Item($r, 6).Value = "Camarillo Aburto Raymundo"
    }
}

# ---------------------------------------------------------------------
# 4) Rescatables: refresh the pending make-up list. Several students who
#    used to be missing BIOLOGIA now have a grade, so they drop off the
#    list; a new set of students (rows 4-10) replace them.
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Rescatables")

# Row 2 keeps its student but the BIOLOGIA docente swaps like on Blancos.
$ws5.Cells.Item(2, 6).Value = "Castro Vasquez Julieta"

# Row 3 is untouched (student's INGLES III row).

# Rows 4-10 become a fresh block of pending make-ups.
$ws5.Cells.Item(4, 1).Value = 20330051920347
$ws5.Cells.Item(4, 2).Value = "LOPEZ"
$ws5.Cells.Item(4, 3).Value = "DE JESUS"
$ws5.Cells.Item(4, 4).Value = "EVELYN"
$ws5.Cells.Item(4, 5).Value = "APLICA LA METODOLOGÍA DE DESARROLLO RÁPIDO DE APLICACIONES CON PROGRAMACIÓN ORIENTADA A EVENTOS"
$ws5.Cells.Item(4, 6).Value = "De Jesús Orduña Sofía del Pilar"
$ws5.Cells.Item(4, 7).Value = -1

$ws5.Cells.Item(5, 1).Value = 20330051920347
$ws5.Cells.Item(5, 2).Value = "LOPEZ"
$ws5.Cells.Item(5, 3).Value = "DE JESUS"
$ws5.Cells.Item(5, 4).Value = "EVELYN"
$ws5.Cells.Item(5, 5).Value = "ÉTICA"
$ws5.Cells.Item(5, 6).Value = "Camarillo Aburto Raymundo"
$ws5.Cells.Item(5, 7).Value = -1

$ws5.Cells.Item(6, 1).Value = 20330051920352
$ws5.Cells.Item(6, 2).Value = "SANCHEZ"
$ws5.Cells.Item(6, 3).Value = "RODRIGUEZ"
$ws5.Cells.Item(6, 4).Value = "EMILIO"
$ws5.Cells.Item(6, 5).Value = "BIOLOGÍA"
$ws5.Cells.Item(6, 6).Value = "Castro Vasquez Julieta"
$ws5.Cells.Item(6, 7).Value = -1

$ws5.Cells.Item(7, 1).Value = 20330051920352
$ws5.Cells.Item(7, 2).Value = "SANCHEZ"
$ws5.Cells.Item(7, 3).Value = "RODRIGUEZ"
$ws5.Cells.Item(7, 4).Value = "EMILIO"
$ws5.Cells.Item(7, 5).Value = "ÉTICA"
$ws5.Cells.Item(7, 6).Value = "Camarillo Aburto Raymundo"
$ws5.Cells.Item(7, 7).Value = -1

$ws5.Cells.Item(8, 1).Value = 20330051920357
$ws5.Cells.Item(8, 2).Value = "XILCAHUA"
$ws5.Cells.Item(8, 3).Value = "TLAXCALA"
$ws5.Cells.Item(8, 4).Value = "LUIS ANGEL"
$ws5.Cells.Item(8, 5).Value = "BIOLOGÍA"
$ws5.Cells.Item(8, 6).Value = "Castro Vasquez Julieta"
$ws5.Cells.Item(8, 7).Value = -1

$ws5.Cells.Item(9, 1).Value = 20330051920357
$ws5.Cells.Item(9, 2).Value = "XILCAHUA"
$ws5.Cells.Item(9, 3).Value = "TLAXCALA"
$ws5.Cells.Item(9, 4).Value = "LUIS ANGEL"
$ws5.Cells.Item(9, 5).Value = "ÉTICA"
$ws5.Cells.Item(9, 6).Value = "Camarillo Aburto Raymundo"
$ws5.Cells.Item(9, 7).Value = -1

$ws5.Cells.Item(10, 1).Value = 20330051920381
$ws5.Cells.Item(10, 2).Value = "HERNANDEZ"
$ws5.Cells.Item(10, 3).Value = "SANCHEZ"
$ws5.Cells.Item(10, 4).Value = "EDGAR DANIEL"
$ws5.Cells.Item(10, 5).Value = "ÉTICA"
$ws5.Cells.Item(10, 6).Value = "Camarillo Aburto Raymundo"
$ws5.Cells.Item(10, 7).Value = -1

# Rows 11-15 no longer apply (those students now have a BIOLOGIA grade) -
# drop them so the sheet shrinks back down to A1:G10.
$ws5.Rows("11:15").Delete()
